# Regolamento erFantacalcio — update season references from 2024-2025 to
# 2025-2026 and adjust the season's first/last matchday dates accordingly.

$d = $word.ActiveDocument

# 1) Title: "Regolamento erFantacalcio 2024-2025" -> "... 2025-2026"
$d.Content.Find.Execute("erFantacalcio 2024-2025", $true, $false, $false, $false, $false, $true, 1, $false, "erFantacalcio 2025-2026", 2)

# 2) Calendar paragraph: first matchday date "sabato 14 settembre" -> "venerdì 19 settembre"
$d.Content.Find.Execute("prevista per sabato 14 settembre", $true, $false, $false, $false, $false, $true, 1, $false, "prevista per venerdì 19 settembre", 2)

# 3) Calendar paragraph: last matchday date "il 25 maggio 2025" -> "il 24 maggio 2026"
$d.Content.Find.Execute("prevista per il 25 maggio 2025", $true, $false, $false, $false, $false, $true, 1, $false, "prevista per il 24 maggio 2026", 2)
